# "Generate Report for handoff"
# The dae061a2-... row has been handed off again, so its status moves from
# "Handed back: in sync with en-US" to "Ready for handoff" on every sheet
# that tracks it, and the per-locale "Latest Handoff Datetime" is stamped
# with the handoff time on the zh-cn and de-de detail sheets.

$wb = $excel.ActiveWorkbook

$statusReady = "Ready for handoff"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $statusReady
$overview.Range("C3").Value = $statusReady

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("B3").Value = $statusReady
$zhcn.Range("D3").Value = "2016-01-25 13:42:18"

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("B3").Value = $statusReady
$dede.Range("D3").Value = "2016-01-25 13:42:28"
